$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.428.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.53%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.811.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +4.17%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'316.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.98%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.01%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.5492"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +10.42%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.3860"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +8.77%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("B9").Value = "'Dogecoin"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.07600"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +4.85%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = "'OKB"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'43.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.42%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'1.134"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +7.03%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.03%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'21.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +6.09%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'6.229"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.51%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = "'Chainlink"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'7.333"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +7.13%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = "'WrappedEther"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'1.811.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.59%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'91.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +5.55%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'0.00001075"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.98%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.06478"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.41%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.00%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'17.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +4.24%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'5.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.54%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'28.451.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +5.37%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'11.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.01%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'2.108"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.08%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'20.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +4.45%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("B27").Value = "'LidoDAOToken"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'2.442"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +14.67%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("B28").Value = "'Monero"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'156.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.02%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'2.015.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.32%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'124.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.99%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'1.174"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +11.04%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'0.1038"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +10.24%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'5.765"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +7.20%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'3.640"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.90%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'0.2326"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +16.27%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("B36").Value = "'Hedera"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.06301"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +6.73%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("B37").Value = "'VeChain"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.02335"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +6.42%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("B38").Value = "'FraxShare"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'8.919"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +19.56%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'11.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +5.22%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.6406"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +6.63%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'5.042"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.97%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'1.170"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +5.62%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.9999"
$ws.Range("D43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'1.389"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.43%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'13.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.16%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.6009"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +6.72%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'3.691"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.27%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'124.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.57%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'1.978"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +6.63%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'1.150"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.72%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.06938"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +4.16%  "
$ws.Range("E51").Style = "Normal"
